$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.65726433333334
$ws.Range("H2").Value = 100.971793
$ws.Range("I2").Value = 0.8115737688004754
$ws.Range("J2").Value = 0.8115737688004754
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 1780.816076465936
$ws.Range("R2").Value = 16027.34468819342
$ws.Range("S2").Value = 0.3377295600469867
$ws.Range("T2").Value = 0.3377295600469867
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.65726433333334
$ws.Range("H3").Value = 100.971793
$ws.Range("I3").Value = 0.8115737688004754
$ws.Range("J3").Value = 0.8115737688004754
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 1590.644275733753
$ws.Range("R3").Value = 14315.79848160378
$ws.Range("S3").Value = 0.3016637139198097
$ws.Range("T3").Value = 0.3016637139198097
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.65726433333334
$ws.Range("H4").Value = 100.971793
$ws.Range("I4").Value = 0.8115737688004754
$ws.Range("J4").Value = 0.8115737688004754
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 907.8914893058719
$ws.Range("R4").Value = 8171.023403752846
$ws.Range("S4").Value = 0.172180494833679
$ws.Range("T4").Value = 0.172180494833679
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.56955
$ws.Range("H5").Value = 1.70865
$ws.Range("I5").Value = 0.01373349406661455
$ws.Range("J5").Value = 0.01373349406661455
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 30.13506345335
$ws.Range("R5").Value = 271.21557108015
$ws.Range("S5").Value = 0.005715077405570918
$ws.Range("T5").Value = 0.005715077405570918
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.56955
$ws.Range("H6").Value = 1.70865
$ws.Range("I6").Value = 0.01373349406661455
$ws.Range("J6").Value = 0.01373349406661455
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 26.91696622374999
$ws.Range("R6").Value = 242.25269601375
$ws.Range("S6").Value = 0.005104769257579519
$ws.Range("T6").Value = 0.005104769257579519
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.56955
$ws.Range("H7").Value = 1.70865
$ws.Range("I7").Value = 0.01373349406661455
$ws.Range("J7").Value = 0.01373349406661455
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 15.3633876067
$ws.Range("R7").Value = 138.2704884603
$ws.Range("S7").Value = 0.00291364740346411
$ws.Range("T7").Value = 0.00291364740346411
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.244787666666667
$ws.Range("H8").Value = 21.734363
$ws.Range("I8").Value = 0.1746927371329101
$ws.Range("J8").Value = 0.1746927371329101
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 383.3239154438548
$ws.Range("R8").Value = 3449.915238994693
$ws.Range("S8").Value = 0.07269690510389873
$ws.Range("T8").Value = 0.07269690510389873
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.244787666666667
$ws.Range("H9").Value = 21.734363
$ws.Range("I9").Value = 0.1746927371329101
$ws.Range("J9").Value = 0.1746927371329101
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 342.3890877392805
$ws.Range("R9").Value = 3081.501789653525
$ws.Range("S9").Value = 0.06493366580368932
$ws.Range("T9").Value = 0.06493366580368932
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.244787666666667
$ws.Range("H10").Value = 21.734363
$ws.Range("I10").Value = 0.1746927371329101
$ws.Range("J10").Value = 0.1746927371329101
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 195.4253025217096
$ws.Range("R10").Value = 1758.827722695386
$ws.Range("S10").Value = 0.037062166225322
$ws.Range("T10").Value = 0.037062166225322
